$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Save" header in H1, reusing the same formatting as the
# existing header cells (bold, bordered, centered) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the new data value in H2 (numeric 0), plain/unstyled like the other
# numeric data cells in row 2.
$ws.Range("H2").Value = 0
